$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed threshold values (Max column) while leaving the
# Parameter labels and everything else in the same row positions.
$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 12
$ws.Range("C5").Value = 20

# Widen columns A and C to fit the longer parameter names (column B keeps
# its existing width).
$ws.Columns("A").ColumnWidth = 26.29
$ws.Columns("C").ColumnWidth = 26.57

# Update the selected cell shown when the sheet is reopened.
$ws.Range("C3").Select()
